$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that currently sits right
#    after the title heading -- its text is being relocated to the very
#    end of the document (see step 3 below).
$d.Paragraphs.Item(2).Range.Delete()

# 2. Build a new bold paragraph containing the page title, using a
#    plain paragraph (no pStyle / no run formatting) as the template so
#    the freshly minted paragraph mark doesn't inherit heading or
#    italic formatting from a neighbour.
$templatePara = $d.Paragraphs.Item(3)
$templatePara.Range.InsertParagraphAfter()
$builtPara = $d.Paragraphs.Item(4)
$insertPoint = $d.Range($builtPara.Range.Start, $builtPara.Range.Start)
$insertPoint.InsertAfter("Play Excalibur Unleashed Free Slot Game - Review")
$boldRange = $d.Range($builtPara.Range.Start, $builtPara.Range.Start + 49)
$boldRange.Font.Bold = 1

# Move the paragraph we just built down to sit right before the closing
# ("image prompt") paragraph by cutting it out and pasting it back in
# at the right spot -- this preserves its own run formatting without
# bleeding into the surrounding paragraphs.
$builtPara.Range.Cut()
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$target = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$target.Paste()

# 3. Replace the old "image generation prompt" text on the final
#    (italic) paragraph with the meta-description copy that used to
#    live at the top of the document.
$d.Content.Find.Execute("Please create a feature image for the game " + [char]34 + "Excalibur Unleashed" + [char]34 + " that features a happy Maya warrior with glasses in a cartoon style. Keep in mind the medieval theme of the game and use colors that match the forest setting of the slot machine. The image should be eye-catching and engaging for players to draw them towards the game.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Experience the legend of Arthur with Excalibur Unleashed free slot game. Read our review to learn about its features, graphics, and high win potential.", `
    2)
